$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column D ("image") data for existing rows 2 and 3 ---
$ws.Cells.Item(2, 4).Value = 1234
$ws.Cells.Item(3, 4).Value = 1234

# --- Add new column G ("link") data for existing rows 2 and 3 ---
$ws.Cells.Item(2, 7).Value = "rjqwehrew"
$ws.Cells.Item(3, 7).Value = "rjqwehrew"

# --- Update the date in row 3 (F3) from 2022-05-12 to 2022-05-10 ---
# (use the raw date serial number so the cell keeps its existing date style)
$ws.Cells.Item(3, 6).Value = 44691

# --- Add a brand-new row 4 ---
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Đề án cuộc thi NCKH"
$ws.Cells.Item(4, 3).Value = "Đề án cuộc thi"
$ws.Cells.Item(4, 4).Value = "fadfaf"
$ws.Cells.Item(4, 5).Value = "BanDaoTao"
$ws.Cells.Item(4, 6).Value = 44692
$ws.Cells.Item(4, 7).Value = "rjqwehrew"

# Copy the date style from F2 onto the new F4 cell so it keeps the date number format
$ws.Cells.Item(2, 6).Copy() | Out-Null
$ws.Cells.Item(4, 6).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Column width adjustments (closest achievable values given COM rounding) ---
$ws.Columns.Item(2).ColumnWidth = 30.26
$ws.Columns.Item(3).ColumnWidth = 16.42
$ws.Columns.Item(5).ColumnWidth = 15.26

# --- Update the selection to reflect the new active range ---
$ws.Range("G2:G4").Select() | Out-Null
